# Update CURRENT STATUS sheet with new run values (19-08-25 run).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns: A=Stock B=Date C=Close D=S_EMA1 E=S_EMA2 F=C_EMA1 G=C_EMA2
#          H=S_RSI_P I=S_RSI_L_R J=C_RSI_V K=P_Pos L=C_Pos
# Each row entry: Row, B(date), C(close), F(C_EMA1), G(C_EMA2), J(C_RSI_V), L(C_Pos or $null if unchanged)
$updates = @(
    @(2,  45883, 85.95,  85.35992399872514, 84.47878764942587, 54.85323046071112, $null),
    @(3,  45883, 115.33, 114.5484223139029, 100.4438058193916, 60.42799201266769, $null),
    @(4,  45883, 48.23,  48.64354744749035, 49.46257330879806, 41.54331864621257, $null),
    @(5,  45883, 263.12, 264.1554814769239, 264.2562112600376, 46.2154049024081,  $null),
    @(6,  45883, 90.5,   90.43098969848623, 91.10364371890013, 43.05394628655823, $null),
    @(7,  45883, 168.06, 170.5588053466418, 172.4218543271275, 45.34208870124789, $null),
    @(8,  45883, 21.69,  21.71946078661329, 21.8709076969679,  47.40499544770443, $null),
    @(9,  45883, 24.76,  24.59706781067737, 24.50845710350486, 58.47951187331634, $null),
    @(10, 45883, 58.35,  58.87610909033982, 58.93253868876573, 46.56604679968317, "S"),
    @(11, 45883, 150.22, 148.5902065783553, 145.4583809894505, 55.12602049899276, $null),
    @(12, 45883, 204.89, 203.659196910959,  186.1799675974528, 63.30069169431682, $null),
    @(13, 45883, 379.83, 378.5934666416255, 401.8534261628802, 88.01200520444478, $null),
    @(14, 45883, 29.82,  29.67945214558003, 29.81722592974812, 49.1071985498797,  "S"),
    @(15, 45883, 567,    568.4002264034567, 571.4143335538879, 46.26793764435742, $null),
    @(16, 45883, 3265.4, 3132.922995600512, 2982.59759644368,  61.67310218031128, $null),
    @(17, 45883, 1936.6, 2105.49035784043,  2236.581348381601, 22.05728328164675, $null)
)

foreach ($u in $updates) {
    $row = $u[0]
    $ws.Cells.Item($row, 2).Value = $u[1]   # B - Date
    $ws.Cells.Item($row, 3).Value = $u[2]   # C - Close
    $ws.Cells.Item($row, 6).Value = $u[3]   # F - C_EMA1
    $ws.Cells.Item($row, 7).Value = $u[4]   # G - C_EMA2
    $ws.Cells.Item($row, 10).Value = $u[5]  # J - C_RSI_V
    if ($null -ne $u[6]) {
        $ws.Cells.Item($row, 12).Value = $u[6]  # L - C_Pos
    }
}
